$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 223117
$ws.Range("E2").Value = -4066
$ws.Range("F2").Value = -2917
$ws.Range("G2").Value = -13286
$ws.Range("H2").Value = -9662
$ws.Range("I2").Value = -10550
$ws.Range("J2").Value = 888
$ws.Range("K2").Value = 337755
$ws.Range("L2").Value = 219852
$ws.Range("M2").Value = 117903
$ws.Range("N2").Value = 103410
$ws.Range("O2").Value = 14493
$ws.Range("P2").Value = 15645
$ws.Range("Q2").Value = 19164
$ws.Range("R2").Value = -31710
$ws.Range("S2").Value = 10717
$ws.Range("T2").Value = 28529
$ws.Range("U2").Value = -9365
$ws.Range("V2").Value = 128704
$ws.Range("W2").Value = -1.82
$ws.Range("X2").Value = -4.33
$ws.Range("Y2").Value = -9.550000000000001
$ws.Range("Z2").Value = -2.82
$ws.Range("AA2").Value = 186.47
$ws.Range("AB2").Value = 639.91
$ws.Range("AC2").Value = -4040
$ws.Range("AD2").Value = -7.73
$ws.Range("AE2").Value = 42232
$ws.Range("AF2").Value = 0.74
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 261111808

# Row 3
$ws.Range("D3").Value = 222812
$ws.Range("E3").Value = 12929
$ws.Range("F3").Value = 12929
$ws.Range("G3").Value = 7195
$ws.Range("H3").Value = 6313
$ws.Range("I3").Value = 5530
$ws.Range("J3").Value = 783
$ws.Range("K3").Value = 293412
$ws.Range("L3").Value = 171757
$ws.Range("M3").Value = 121655
$ws.Range("N3").Value = 108451
$ws.Range("O3").Value = 13204
$ws.Range("P3").Value = 15645
$ws.Range("Q3").Value = 42300
$ws.Range("R3").Value = -24019
$ws.Range("S3").Value = -11640
$ws.Range("T3").Value = 31157
$ws.Range("U3").Value = 11142
$ws.Range("V3").Value = 87911
$ws.Range("W3").Value = 5.8
$ws.Range("X3").Value = 2.83
$ws.Range("Y3").Value = 5.22
$ws.Range("Z3").Value = 2
$ws.Range("AA3").Value = 141.18
$ws.Range("AB3").Value = 671.3
$ws.Range("AC3").Value = 2118
$ws.Range("AD3").Value = 13.34
$ws.Range("AE3").Value = 44293
$ws.Range("AF3").Value = 0.64
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 1.77
$ws.Range("AI3").Value = 22.14
$ws.Range("AJ3").Value = 261111808

# Row 4
$ws.Range("D4").Value = 227437
$ws.Range("E4").Value = 14400
$ws.Range("F4").Value = 14400
$ws.Range("G4").Value = 11270
$ws.Range("H4").Value = 7978
$ws.Range("I4").Value = 7111
$ws.Range("J4").Value = 868
$ws.Range("K4").Value = 305877
$ws.Range("L4").Value = 177930
$ws.Range("M4").Value = 127948
$ws.Range("N4").Value = 114419
$ws.Range("O4").Value = 13528
$ws.Range("P4").Value = 15645
$ws.Range("Q4").Value = 47708
$ws.Range("R4").Value = -34850
$ws.Range("S4").Value = -9433
$ws.Range("T4").Value = 27643
$ws.Range("U4").Value = 20064
$ws.Range("V4").Value = 83015
$ws.Range("W4").Value = 6.33
$ws.Range("X4").Value = 3.51
$ws.Range("Y4").Value = 6.38
$ws.Range("Z4").Value = 2.66
$ws.Range("AA4").Value = 139.06
$ws.Range("AB4").Value = 709.33
$ws.Range("AC4").Value = 2723
$ws.Range("AD4").Value = 10.8
$ws.Range("AE4").Value = 46707
$ws.Range("AF4").Value = 0.63
$ws.Range("AG4").Value = 800
$ws.Range("AH4").Value = 2.72
$ws.Range("AI4").Value = 27.56
$ws.Range("AJ4").Value = 261111808

# Row 5
$ws.Range("D5").Value = 233873
$ws.Range("E5").Value = 13753
$ws.Range("F5").Value = 13753
$ws.Range("G5").Value = 8370
$ws.Range("H5").Value = 5615
$ws.Range("I5").Value = 4767
$ws.Range("J5").Value = 848
$ws.Range("K5").Value = 297309
$ws.Range("L5").Value = 165203
$ws.Range("M5").Value = 132106
$ws.Range("N5").Value = 118188
$ws.Range("O5").Value = 13918
$ws.Range("P5").Value = 15645
$ws.Range("Q5").Value = 38777
$ws.Range("R5").Value = -34832
$ws.Range("S5").Value = -13635
$ws.Range("T5").Value = 24422
$ws.Range("U5").Value = 14355
$ws.Range("V5").Value = 67764
$ws.Range("W5").Value = 5.88
$ws.Range("X5").Value = 2.4
$ws.Range("Y5").Value = 4.1
$ws.Range("Z5").Value = 1.86
$ws.Range("AA5").Value = 125.05
$ws.Range("AB5").Value = 730.5599999999999
$ws.Range("AC5").Value = 1826
$ws.Range("AD5").Value = 16.57
$ws.Range("AE5").Value = 48221
$ws.Range("AF5").Value = 0.63
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 3.31
$ws.Range("AI5").Value = 51.41
$ws.Range("AJ5").Value = 261111808

# Row 6
$ws.Range("D6").Value = 234601
$ws.Range("E6").Value = 12615
$ws.Range("F6").Value = 12615
$ws.Range("G6").Value = 10907
$ws.Range("H6").Value = 7623
$ws.Range("I6").Value = 6885
$ws.Range("K6").Value = 321888
$ws.Range("L6").Value = 174576
$ws.Range("M6").Value = 147313
$ws.Range("N6").Value = 132027
$ws.Range("P6").Value = 15645
$ws.Range("Q6").Value = 40105
$ws.Range("R6").Value = -27041
$ws.Range("S6").Value = -5317
$ws.Range("T6").Value = 22609
$ws.Range("U6").Value = 17496
$ws.Range("V6").Value = 67554
$ws.Range("W6").Value = 5.38
$ws.Range("X6").Value = 3.25
$ws.Range("Y6").Value = 5.5
$ws.Range("Z6").Value = 2.46
$ws.Range("AA6").Value = 118.51
$ws.Range("AB6").Value = 816.1799999999999
$ws.Range("AC6").Value = 2637
$ws.Range("AD6").Value = 11.3
$ws.Range("AE6").Value = 53857
$ws.Range("AF6").Value = 0.55
$ws.Range("AI6").Value = 39.17
$ws.Range("AJ6").Value = 261111808

# Row 7
$ws.Range("D7").Value = 242775
$ws.Range("E7").Value = 11764
$ws.Range("G7").Value = 10506
$ws.Range("H7").Value = 7530
$ws.Range("I7").Value = 6785
$ws.Range("K7").Value = 329217
$ws.Range("L7").Value = 177128
$ws.Range("M7").Value = 152088
$ws.Range("N7").Value = 136278
$ws.Range("P7").Value = 15643
$ws.Range("Q7").Value = 43577
$ws.Range("R7").Value = -40394
$ws.Range("S7").Value = -6956
$ws.Range("T7").Value = 34541
$ws.Range("U7").Value = 5258
$ws.Range("W7").Value = 4.85
$ws.Range("X7").Value = 3.1
$ws.Range("Y7").Value = 5.06
$ws.Range("Z7").Value = 2.31
$ws.Range("AA7").Value = 116.46
$ws.Range("AC7").Value = 2598
$ws.Range("AD7").Value = 9.76
$ws.Range("AE7").Value = 55569
$ws.Range("AF7").Value = 0.46
$ws.Range("AG7").Value = 1090
$ws.Range("AH7").Value = 4.3
$ws.Range("AI7").Value = 41.97

# Row 8
$ws.Range("D8").Value = 250319
$ws.Range("E8").Value = 12699
$ws.Range("G8").Value = 11402
$ws.Range("H8").Value = 8498
$ws.Range("I8").Value = 7616
$ws.Range("K8").Value = 335047
$ws.Range("L8").Value = 177527
$ws.Range("M8").Value = 157520
$ws.Range("N8").Value = 141017
$ws.Range("P8").Value = 15643
$ws.Range("Q8").Value = 43616
$ws.Range("R8").Value = -39912
$ws.Range("S8").Value = -3153
$ws.Range("T8").Value = 31244
$ws.Range("U8").Value = 10369
$ws.Range("W8").Value = 5.07
$ws.Range("X8").Value = 3.4
$ws.Range("Y8").Value = 5.49
$ws.Range("Z8").Value = 2.56
$ws.Range("AA8").Value = 112.7
$ws.Range("AC8").Value = 2917
$ws.Range("AD8").Value = 8.69
$ws.Range("AE8").Value = 57501
$ws.Range("AF8").Value = 0.44
$ws.Range("AG8").Value = 1114
$ws.Range("AH8").Value = 4.4
$ws.Range("AI8").Value = 38.2

# Row 9
$ws.Range("D9").Value = 258145
$ws.Range("E9").Value = 14250
$ws.Range("G9").Value = 12932
$ws.Range("H9").Value = 9559
$ws.Range("I9").Value = 8651
$ws.Range("K9").Value = 344465
$ws.Range("L9").Value = 180844
$ws.Range("M9").Value = 163622
$ws.Range("N9").Value = 146490
$ws.Range("P9").Value = 15643
$ws.Range("Q9").Value = 45298
$ws.Range("R9").Value = -38795
$ws.Range("S9").Value = -2665
$ws.Range("T9").Value = 29124
$ws.Range("U9").Value = 14114
$ws.Range("W9").Value = 5.52
$ws.Range("X9").Value = 3.7
$ws.Range("Y9").Value = 6.02
$ws.Range("Z9").Value = 2.81
$ws.Range("AA9").Value = 110.53
$ws.Range("AC9").Value = 3313
$ws.Range("AD9").Value = 7.65
$ws.Range("AE9").Value = 59733
$ws.Range("AF9").Value = 0.42
$ws.Range("AG9").Value = 1158
$ws.Range("AH9").Value = 4.57
$ws.Range("AI9").Value = 34.95

# Remove cells that no longer exist in the target (structural deletion)
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
